$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "325÷5=" "125÷7="
Replace-Text "789÷3=" "888÷7="
Replace-Text "952÷2=" "359÷9="
Replace-Text "532÷2=" "601÷7="
Replace-Text "642÷7=" "968÷4="
Replace-Text "133÷3=" "566÷5="
Replace-Text "116÷2=" "684÷5="
Replace-Text "655÷2=" "176÷9="
Replace-Text "187÷4=" "405÷3="
Replace-Text "834÷8=" "701÷9="
Replace-Text "943÷9=" "237÷2="
Replace-Text "189÷7=" "966÷3="
Replace-Text "884÷8=" "221÷5="
Replace-Text "973÷7=" "711÷3="
Replace-Text "186÷2=" "396÷2="
Replace-Text "662÷5=" "483÷4="
Replace-Text "236÷3=" "357÷4="
Replace-Text "231÷5=" "318÷4="
Replace-Text "605÷4=" "404÷5="
Replace-Text "385÷4=" "622÷5="
Replace-Text "802÷8=" "175÷8="
Replace-Text "145÷8=" "477÷2="
Replace-Text "304÷2=" "832÷2="
Replace-Text "493÷4=" "910÷4="
Replace-Text "376÷9=" "764÷4="
